$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the "Se connecter au site web" use case to "Se connecter" (row 3, column A).
$ws.Range("A3").Value = "Se connecter"

# 2. Extend the acceptance criteria of "Accéder à une conversation" (row 9, column C) so it
#    also covers seeing whether a sent message has been read and seeing the send time of
#    messages (this absorbs the acceptance criteria of the deleted "Visualiser la lecture
#    d'un message" user story below).
$c9 = $ws.Range("C9")
$c9.Value = "Etant donné que je suis connecté et sur l'interface de chat" + "`n" + "Lorsque je clique sur un onglet existant d'une conversation" + "`n" + "Alors je peux accéder à la conversation, voir si un message que j'ai envoyé a été lu et voir l'heure d'envoi des messages"

# Re-bold the three keyword runs ("Etant donné que", "Lorsque", "Alors") that introduce each
# part of the Gherkin-style acceptance criteria, matching the original rich-text formatting.
$c9.Characters(1, 15).Font.Bold = $true
$c9.Characters(61, 7).Font.Bold = $true
$c9.Characters(121, 5).Font.Bold = $true

# 3. Remove the "Visualiser la lecture d'un message" user story entirely (use case, user
#    story and acceptance-criteria cells) - it is now covered by the extended acceptance
#    criteria of "Accéder à une conversation" above.
$ws.Rows.Item(10).Delete()
